$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new blank row at position 7 (splits the old "studydesign.csv" row
#     into two rows: "dose per bodyweight" and "dose per surface area") ---
$ws.Rows.Item(7).Insert()

# --- Row 3: "Check if output structure is correctly set" test replaced with
#     the new "Find inconsitencies in output cvs" test ---
$ws.Range("C3").Value = "Find inconsitencies in output cvs"

$d3 = @"
Breaks with logfile message
ERROR: Outputpath "Organism|PeripheralVenousBlood|Hydroxy_Itraconazole|Plasma (Peripheral Venous Blood) with Typo" could not be found in model
ERROR: For unit "µmol/l", there is no common dimension with display unit "cm"
ERROR: unit "typo" for seems to be no default OSPSuite unit
ERROR: unit "typo2" for seems to be no default OSPSuite unit
ERROR: For unit "µmol/l", there is no common dimension with display unit "h"
"@
$ws.Range("D3").Value = $d3

# --- Row 4: single-application test - text stays the same but font color
#     changes from red to black (automatic), and a new "expected outcome"
#     cell is added ---
$ws.Range("C4").Font.Color = 0

$d4 = @"
Two csv files are generated: 
simulations/SingleIvBolus-Results.csv
simulations/SingleIvBolus-PK-Analyses.csv
"@
$ws.Range("D4").Value = $d4

# --- Row 5: multi-application test - fix typo "caluclation" -> "calculation" ---
$ws.Range("C5").Value = "Population simulation and PK Parameter calculation of a multi application"

# --- Row 6: first half of the studydesign.csv split - "dose per bodyweigth" ---
$ws.Range("C6").Value = "Populations simulation with studydesign.csv dose per bodyweigth"
$ws.Range("C6").Font.Color = 0

$d6 = @"
Two csv files are generated: 
simulations/SingleIvBolus-Results.csv
simulations/SingleIvBolus-PK-Analyses.csv
"@
$ws.Range("D6").Value = $d6
$ws.Rows.Item(6).RowHeight = 41.4

# --- Row 7: second half of the studydesign.csv split - "dose per surface area" ---
$ws.Range("B7").Value = 6

$ws.Range("C7").Value = "Populations simulation with studydesign.csv dose per surface area"
$ws.Range("C7").Font.Color = 0

$d7 = @"
Two csv files are generated: 
simulations/SingleIvBolus-Results.csv
simulations/SingleIvBolus-PK-Analyses.csv
"@
$ws.Range("D7").Value = $d7

$ws.Range("E7").Value = "7.2_BSA_Example"

$ws.Rows.Item(7).RowHeight = 41.4

# --- Rows 8-12 (previously 7-11): renumber column B, the row text shifted
#     down automatically with the earlier Insert() ---
$ws.Range("B8").Value = 7
$ws.Range("B9").Value = 8
$ws.Range("B10").Value = 9
$ws.Range("B11").Value = 10
$ws.Range("B12").Clear()

# --- Sheet view: scroll down a bit, move selection to D11 ---
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("D11").Select()

$wb.Save()
